# "fixed null values and added sorting"
# Replace bogus/placeholder E-column values across the dev1..dev10 sheets
# with corrected figures, and update each sheet's remembered selection to
# reflect where the user last clicked while reviewing the data. dev10 ends
# up as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# dev1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev1")
$ws.Range("E3").Value = 13434
$ws.Range("E4").Value = 14343
$ws.Range("E5").Value = 15342
$ws.Range("E6").Value = 16342
$ws.Range("E7").Value = 17121
$ws.Range("E8").Value = 17900
$ws.Range("E9").Value = 18123
$ws.Range("E10").Value = 19234
$ws.Range("E11").Value = 20123
$ws.Range("E12").Value = 22234
$ws.Range("E13").Value = 23123
$ws.Range("E14").Value = 23999
$ws.Range("E15").Value = 24990
$ws.Range("E16").Value = 26000
$ws.Range("E16").Select() | Out-Null

# ---------------------------------------------------------------------
# dev2
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev2")
$ws.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------
# dev3 - untouched
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# dev4 - selection cleared back to the default (top-left cell)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev4")
$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# dev5 - no longer the active tab (handled by activating dev10 below);
# its own in-sheet selection (F10) is unchanged.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# dev6
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev6")
$ws.Range("E3").Value = 33000
$ws.Range("E4").Value = 34000
$ws.Range("E5").Value = 35121
$ws.Range("E6").Value = 36123
$ws.Range("E7").Value = 36812
$ws.Range("E8").Value = 39123
$ws.Range("E9").Value = 40000
$ws.Range("E10").Value = 41221
$ws.Range("E11").Value = 42123
$ws.Range("E11").Select() | Out-Null

# ---------------------------------------------------------------------
# dev7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev7")
$ws.Range("E4").Value = 6234
$ws.Range("E5").Value = 7533
$ws.Range("E6").Value = 8123
$ws.Range("E7").Value = 9145
$ws.Range("E8").Value = 10654
$ws.Range("E9").Value = 12451
$ws.Range("E10").Select() | Out-Null

# ---------------------------------------------------------------------
# dev8
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev8")
$ws.Range("E3").Value = 13213
$ws.Range("E4").Value = 14522
$ws.Range("E5").Value = 15532
$ws.Range("E6").Value = 16093
$ws.Range("E7").Value = 17342
$ws.Range("E8").Value = 18340
$ws.Range("E10").Value = 11234
$ws.Range("E11").Value = 13234
$ws.Range("E12").Value = 14542
$ws.Range("E13").Value = 17934
$ws.Range("E14").Value = 19324
$ws.Range("E15").Value = 20000
$ws.Range("E16").Value = 21453
$ws.Range("E16").Select() | Out-Null

# ---------------------------------------------------------------------
# dev9
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev9")
$ws.Range("E3").Value = 7343
$ws.Range("E4").Value = 8231
$ws.Range("E5").Value = 9422
$ws.Range("E6").Value = 10234
$ws.Range("E7").Value = 12430
$ws.Range("E8").Value = 13401
$ws.Range("E9").Value = 14933
$ws.Range("E10").Value = 16123
$ws.Range("E11").Value = 17390
$ws.Range("E12").Value = 18231
$ws.Range("E13").Value = 19234
$ws.Range("E14").Value = 20542
$ws.Range("E15").Value = 22034
$ws.Range("E16").Value = 23400
$ws.Range("E17").Value = 24013
$ws.Range("E18").Value = 27193
$ws.Range("E19").Value = 29012
$ws.Range("E19").Select() | Out-Null

# ---------------------------------------------------------------------
# dev10 - becomes the active tab
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dev10")
$ws.Range("E2").Value = 1231
$ws.Range("E3").Value = 4231
$ws.Range("E4").Value = 5231
$ws.Range("E5").Value = 7231
$ws.Range("E6").Value = 8123
$ws.Activate() | Out-Null
$ws.Range("E6").Select() | Out-Null
